$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.420.17"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "1.880.39"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7173"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07947"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3148"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08144"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.84%  "
$ws.Range("D12").Value = "1.894.51"
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "94.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.01%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.244"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7090"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.393"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008436"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.87%  "
$ws.Range("D18").Value = "29.434.71"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "253.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.75%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.144.28"
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.732"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.075"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.510"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.415"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.292"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.226"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05321"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.947"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7581"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.14%  "
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.700"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01899"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.35%  "
$ws.Range("D39").Value = "1.275.53"
$ws.Range("E39").Value = "  +2.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.759"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.410"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "112.93"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.74%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9067"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.59%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "74.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.70%  "
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("D47").Value = "2.038.91"
$ws.Range("E47").Value = "  +1.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.810"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5207"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.504"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4353"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.36%  "
